$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 991, shifting existing rows 991-1038 down to 992-1039.
$ws.Rows.Item(991).Insert()

# Populate the newly inserted row 991 with the new data point.
$ws.Range("A991").Value = 10
$ws.Range("B991").Value = "Vega Modelo de Temuco"
$ws.Range("C991").Value = "La Araucanía"
$ws.Range("D991").Value = 45041
$ws.Range("E991").Value = 9
$ws.Range("F991").Value = 100112003
$ws.Range("G991").Value = "Ajo"
$ws.Range("H991").Value = "Chino"
$ws.Range("I991").Value = "Primera"
$ws.Range("J991").Value = 500
$ws.Range("K991").Value = 18000
$ws.Range("L991").Value = 18000
$ws.Range("M991").Value = 18000
$ws.Range("N991").Value = "$/caja 10 kilos"
$ws.Range("O991").Value = "China"
$ws.Range("P991").Value = 1800
$ws.Range("Q991").Value = 10
$ws.Range("R991").Value = "Hortaliza"
